{"js": "// Update the date line and the 25 division problems to the new values\n// from the commit diff. Each old text string in this document is\n// unique, so a body-wide search/replace of the exact run text is safe\n// and keeps the original run formatting (font, size, etc.) intact.\nconst replacements = [\n  [\"2024-04-27 Saturday\", \"2024-04-28 Sunday\"],\n  [\"18\u00f77=\", \"69\u00f73=\"],\n  [\"54\u00f76=\", \"92\u00f74=\"],\n  [\"76\u00f73=\", \"77\u00f79=\"],\n  [\"45\u00f73=\", \"62\u00f76=\"],\n  [\"78\u00f76=\", \"28\u00f72=\"],\n  [\"65\u00f78=\", \"59\u00f78=\"],\n  [\"11\u00f76=\", \"41\u00f75=\"],\n  [\"96\u00f72=\", \"52\u00f73=\"],\n  [\"69\u00f72=\", \"27\u00f73=\"],\n  [\"24\u00f77=\", \"69\u00f77=\"],\n  [\"76\u00f75=\", \"32\u00f76=\"],\n  [\"20\u00f72=\", \"33\u00f78=\"],\n  [\"24\u00f74=\", \"43\u00f76=\"],\n  [\"40\u00f74=\", \"33\u00f72=\"],\n  [\"14\u00f76=\", \"65\u00f79=\"],\n  [\"59\u00f76=\", \"87\u00f74=\"],\n  [\"17\u00f76=\", \"63\u00f74=\"],\n  [\"38\u00f78=\", \"69\u00f76=\"],\n  [\"87\u00f79=\", \"93\u00f78=\"],\n  [\"14\u00f73=\", \"65\u00f77=\"],\n  [\"27\u00f75=\", \"25\u00f78=\"],\n  [\"97\u00f74=\", \"81\u00f79=\"],\n  [\"71\u00f77=\", \"86\u00f76=\"],\n  [\"20\u00f79=\", \"48\u00f72=\"],\n  [\"71\u00f79=\", \"97\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division problems to the new values\n# from the commit diff. Each old text string in this document is\n# unique, so Find/Replace of the exact text (match case, no wildcards)\n# is safe and preserves the original run formatting (font, size, etc.).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-27 Saturday\", \"2024-04-28 Sunday\"),\n    @(\"18\u00f77=\", \"69\u00f73=\"),\n    @(\"54\u00f76=\", \"92\u00f74=\"),\n    @(\"76\u00f73=\", \"77\u00f79=\"),\n    @(\"45\u00f73=\", \"62\u00f76=\"),\n    @(\"78\u00f76=\", \"28\u00f72=\"),\n    @(\"65\u00f78=\", \"59\u00f78=\"),\n    @(\"11\u00f76=\", \"41\u00f75=\"),\n    @(\"96\u00f72=\", \"52\u00f73=\"),\n    @(\"69\u00f72=\", \"27\u00f73=\"),\n    @(\"24\u00f77=\", \"69\u00f77=\"),\n    @(\"76\u00f75=\", \"32\u00f76=\"),\n    @(\"20\u00f72=\", \"33\u00f78=\"),\n    @(\"24\u00f74=\", \"43\u00f76=\"),\n    @(\"40\u00f74=\", \"33\u00f72=\"),\n    @(\"14\u00f76=\", \"65\u00f79=\"),\n    @(\"59\u00f76=\", \"87\u00f74=\"),\n    @(\"17\u00f76=\", \"63\u00f74=\"),\n    @(\"38\u00f78=\", \"69\u00f76=\"),\n    @(\"87\u00f79=\", \"93\u00f78=\"),\n    @(\"14\u00f73=\", \"65\u00f77=\"),\n    @(\"27\u00f75=\", \"25\u00f78=\"),\n    @(\"97\u00f74=\", \"81\u00f79=\"),\n    @(\"71\u00f77=\", \"86\u00f76=\"),\n    @(\"20\u00f79=\", \"48\u00f72=\"),\n    @(\"71\u00f79=\", \"97\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,\n        $false, $true, $false, $false, $false,\n        $true, 1, $false,\n        $newText, 2\n    )\n}\n"}
